$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Topic/Parent data (columns A and B), rows 2-37 ---
$ws.Range("A2").Value = "Call"
$ws.Range("A3").Value = "Launch"
$ws.Range("A4").Value = "Malfunction"; $ws.Range("B4").Value = "Launch"
$ws.Range("A5").Value = "Others"; $ws.Range("B5").Value = "Launch"
$ws.Range("A6").Value = "Confront"; $ws.Range("B6").Value = "Launch"
$ws.Range("A7").Value = "Forget"; $ws.Range("B7").Value = "Launch"
$ws.Range("A8").Value = "Normality"
$ws.Range("A9").Value = "Calamity"
$ws.Range("A10").Value = "War"; $ws.Range("B10").Value = "Calamity"
$ws.Range("A11").Value = "Defcon"; $ws.Range("B11").Value = "War"
$ws.Range("A12").Value = "Nukes"
$ws.Range("A13").Value = "Development"; $ws.Range("B13").Value = "Nukes"
$ws.Range("A14").Value = "Destruction"; $ws.Range("B14").Value = "Nukes"
$ws.Range("A15").Value = "Civilians"
$ws.Range("A16").Value = "Authority"
$ws.Range("A17").Value = "Russians"
$ws.Range("A18").Value = "Duty"
$ws.Range("A19").Value = "Complex"
$ws.Range("A20").Value = "Justice"
$ws.Range("A21").Value = "Civilization"
$ws.Range("A22").Value = "Eroticism"
$ws.Range("A23").Value = "History"
$ws.Range("A24").Value = "Japan"; $ws.Range("B24").Value = "History"
$ws.Range("A25").Value = "Indians"; $ws.Range("B25").Value = "History"
$ws.Range("A26").Value = "Scientists"; $ws.Range("B26").Value = "History"
$ws.Range("A27").Value = "Doomsday"
$ws.Range("A28").Value = "Jack"
$ws.Range("A29").Value = "Nature"
$ws.Range("A30").Value = "Man"; $ws.Range("B30").Value = "Nature"
$ws.Range("A31").Value = "Dreamers"; $ws.Range("B31").Value = "Man"
$ws.Range("A32").Value = "Death"; $ws.Range("B32").Value = "Man"
$ws.Range("A33").Value = "Time"; $ws.Range("B33").Value = "Nature"
$ws.Range("A34").Value = "Land"; $ws.Range("B34").Value = "Nature"
$ws.Range("A35").Value = "God"; $ws.Range("B35").Value = "Nature"
$ws.Range("A36").Value = "Cycles"; $ws.Range("B36").Value = "Time"
$ws.Range("A37").Value = "Childhood"; $ws.Range("B37").Value = "Time"

# --- Bold styling for column A topic names (font 16, reused bold/theme1) ---
$ws.Range("A2:A5").Font.Bold = $true
$ws.Range("A8:A22").Font.Bold = $true
$ws.Range("A24:A28").Font.Bold = $true
$ws.Range("A30:A37").Font.Bold = $true

# --- Bold + italic styling for "History" and "Nature" (top-level group headers) ---
$ws.Range("A23").Font.Bold = $true
$ws.Range("A23").Font.Italic = $true
$ws.Range("A29").Font.Bold = $true
$ws.Range("A29").Font.Italic = $true

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection matching final cursor position ---
[void]$ws.Range("C36").Select()
